# Update the "Run_Script" sheet so the dev-environment base URLs are
# repointed from the "development1" host to the new "codetest" host,
# and move the active selection to B21 (matching the author's new
# working cell for the 6-1-24 commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Base URL cells (hyperlink-styled text, targets unchanged, only the
# displayed text changes to the new codetest.advantageclub.co domain)
$ws.Range("B2").Value  = "https://codetest.advantageclub.co"
$ws.Range("B3").Value  = "https://codetest.advantageclub.co/in/rewards/home"
$ws.Range("B21").Value = "https://codetest.advantageclub.co"

# Move/save the active selection to B21
$ws.Range("B21").Select() | Out-Null
